$wb = $excel.ActiveWorkbook

# --- Locate existing sheets -------------------------------------------------
# Current order: 2020-Q4, 2021-Q1, 2021-Q2, 总计
$refSheet = $wb.Worksheets.Item(3)   # "2021-Q2" - used as a style donor (header/col-A style)

# --- Create the new "2022-Q1" sheet ----------------------------------------
# Copy "总计" so the new sheet inherits identical sheet-level structure
# (sheetPr/outlinePr/pageMargins/etc.), then wipe its cell content and
# rebuild it with the 2022-Q1 fund data. NOTE: after Copy() the handle we
# call Copy on ends up referring to the *new* sheet, not the original, so
# re-fetch sheets by name/index afterwards instead of reusing old handles.
$wb.Worksheets.Item("总计").Copy($wb.Worksheets.Item("总计"))

$newQ1 = $wb.Worksheets.Item(4)
$newQ1.Name = "2022-Q1"
$newQ1.Cells.Clear()

# Re-use the exact header / first-column style from an existing quarter sheet.
$refSheet.Range("B1:H1").Copy()
$newQ1.Range("B1:H1").PasteSpecial(-4122)
$refSheet.Range("A2").Copy()
$newQ1.Range("A2").PasteSpecial(-4122)

$newQ1.Range("B1").Value = "基金代码"
$newQ1.Range("C1").Value = "基金名称"
$newQ1.Range("D1").Value = "基金规模"
$newQ1.Range("E1").Value = "股票总仓位"
$newQ1.Range("F1").Value = "仓位占比"
$newQ1.Range("G1").Value = "持有市值(亿元)"
$newQ1.Range("H1").Value = "仓位排名"

$newQ1.Range("A2").Value = 0

$newQ1.Range("B2").NumberFormat = "@"
$newQ1.Range("B2").Value = "159792"
$newQ1.Range("C2").Value = "富国中证港股通互联网ETF"
$newQ1.Range("D2").NumberFormat = "@"
$newQ1.Range("D2").Value = "2.76"
$newQ1.Range("E2").NumberFormat = "@"
$newQ1.Range("E2").Value = "99.00"
$newQ1.Range("F2").NumberFormat = "@"
$newQ1.Range("F2").Value = "2.48"
$newQ1.Range("G2").NumberFormat = "@"
$newQ1.Range("G2").Value = "0.0684"
$newQ1.Range("H2").Value = 10

# --- Update the "总计" sheet -------------------------------------------------
# Insert a new row 2 for 2022-Q1 and push the existing quarters down.
$total = $wb.Worksheets.Item("总计")

$total.Rows(2).Insert()

# Row 2 index cell (column A) should carry the same style as the other index
# cells in column A; copy it from row 3 (the row that used to be row 2).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.07000000000000001

# Renumber the (0-based) index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
